# Refresh the cryptos list (prices / 1h volume %, plus a few rank swaps)
# with the latest scrape. Every cell below is stored as plain text in the
# source workbook (t="inlineStr"), including values that look numeric, e.g.
# "63.147.17" / "5.28" / "1.00" - those are price strings with thousands-dot
# grouping, not real numbers. Writing them through .Value would make Excel
# auto-convert anything that parses as a number, so for those we briefly flip
# the cell to text format, assign the literal string, then restore the default
# "Normal" style so no stray formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# Each triple is (row, column, new cell text).
$updates = @(
    @(2, 4, '63.147.17'),
    @(2, 5, '  +2.12%  '),
    @(3, 4, '2.461.24'),
    @(3, 5, '  +1.94%  '),
    @(4, 5, '  -0.25%  '),
    @(5, 4, '577.65'),
    @(5, 5, '  +1.47%  '),
    @(6, 4, '146.60'),
    @(6, 5, '  +1.99%  '),
    @(7, 5, '  +0.19%  '),
    @(8, 4, '0.542'),
    @(8, 5, '  +0.92%  '),
    @(9, 4, '2.461.04'),
    @(9, 5, '  +1.32%  '),
    @(10, 4, '0.112'),
    @(10, 5, '  +2.34%  '),
    @(11, 5, '  +1.59%  '),
    @(12, 4, '5.28'),
    @(12, 5, '  +1.16%  '),
    @(13, 5, '  +2.11%  '),
    @(14, 5, '  +9.81%  '),
    @(15, 5, '  +2.09%  '),
    @(16, 4, '2.909.09'),
    @(16, 5, '  +2.55%  '),
    @(17, 4, '63.051.63'),
    @(17, 5, '  +2.19%  '),
    @(18, 4, '2.475.00'),
    @(18, 5, '  +2.02%  '),
    @(19, 4, '7.94'),
    @(19, 5, '  -0.96%  '),
    @(20, 4, '11.12'),
    @(20, 5, '  +3.80%  '),
    @(21, 4, '330.53'),
    @(21, 5, '  +1.72%  '),
    @(22, 4, '2.25'),
    @(22, 5, '  +9.64%  '),
    @(23, 5, '  +1.08%  '),
    @(24, 5, '  -0.04%  '),
    @(25, 4, '66.59'),
    @(25, 5, '  +2.16%  '),
    @(26, 4, '667.88'),
    @(26, 5, '  +7.00%  '),
    @(27, 4, '9.06'),
    @(27, 5, '  +8.04%  '),
    @(28, 4, '1.07'),
    @(28, 5, '  +7.39%  '),
    @(29, 5, '  +4.67%  '),
    @(31, 5, '  +2.70%  '),
    @(32, 5, '  +1.27%  '),
    @(33, 4, '1.89'),
    @(33, 5, '  +4.21%  '),
    @(34, 5, '  +3.02%  '),
    @(35, 4, '1.56'),
    @(35, 5, '  +5.54%  '),
    @(36, 4, '0.999'),
    @(36, 5, '  +0.22%  '),
    @(37, 4, '4.81'),
    @(37, 5, '  +3.72%  '),
    @(38, 4, '5.54'),
    @(38, 5, '  +3.32%  '),
    @(39, 2, 'PolygonEcosystemToken'),
    @(39, 3, 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'),
    @(39, 4, '0.374'),
    @(39, 5, '  +0.65%  '),
    @(40, 2, 'Monero'),
    @(40, 3, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @(40, 4, '152.73'),
    @(40, 5, '  -0.31%  '),
    @(41, 5, '  +2.32%  '),
    @(42, 5, '  +6.60%  '),
    @(43, 2, 'Stacks'),
    @(43, 3, 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'),
    @(43, 4, '1.78'),
    @(43, 5, '  +3.32%  '),
    @(44, 2, 'OKB'),
    @(44, 3, 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
    @(44, 4, '42.41'),
    @(44, 5, '  +0.88%  '),
    @(45, 2, 'USDe'),
    @(45, 3, 'https://coinranking.com/coin/exbfr2U-0+usde-usde'),
    @(45, 4, '1.00'),
    @(45, 5, '  +0.11%  '),
    @(46, 2, 'BabyDogeCoin'),
    @(46, 3, 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'),
    @(46, 4, '0.0₆0306'),
    @(46, 5, '  +7.36%  '),
    @(47, 4, '15.13'),
    @(47, 5, '  +27.57%  '),
    @(48, 4, '146.53'),
    @(48, 5, '  +2.43%  '),
    @(49, 5, '  +2.45%  '),
    @(50, 4, '20.80'),
    @(50, 5, '  +3.98%  '),
    @(51, 5, '  +1.73%  ')
)

foreach ($u in $updates) {
    Set-TextCell $u[0] $u[1] $u[2]
}
